# Updates price / volume / (and a few swapped rank rows) for the cryptos sheet,
# matching the GitHub Actions data refresh described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (never let Excel auto-convert
# numeric-looking strings like '1.00' or '599.75' into real numbers),
# then strip the forced-text style back to Normal so no new cell
# formatting is introduced.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = 'Normal'
}

# Row 2
Set-TextValue 'D2' '67.990.52'
Set-TextValue 'E2' '  -2.83%  '
# Row 3
Set-TextValue 'D3' '3.842.64'
Set-TextValue 'E3' '  -2.58%  '
# Row 4
Set-TextValue 'E4' '  +0.34%  '
# Row 5
Set-TextValue 'D5' '599.75'
Set-TextValue 'E5' '  -1.74%  '
# Row 6
Set-TextValue 'D6' '167.83'
Set-TextValue 'E6' '  -1.27%  '
# Row 7
Set-TextValue 'D7' '3.844.41'
Set-TextValue 'E7' '  -2.50%  '
# Row 8
Set-TextValue 'E8' '  +0.03%  '
# Row 9
Set-TextValue 'E9' '  -1.51%  '
# Row 10
Set-TextValue 'E10' '  -5.24%  '
# Row 11
Set-TextValue 'D11' '6.44'
Set-TextValue 'E11' '  -0.73%  '
# Row 12
Set-TextValue 'D12' '0.456'
Set-TextValue 'E12' '  -2.78%  '
# Row 13
Set-TextValue 'D13' '0.0000257'
Set-TextValue 'E13' '  -0.25%  '
# Row 14
Set-TextValue 'D14' '36.86'
Set-TextValue 'E14' '  -3.15%  '
# Row 15
Set-TextValue 'D15' '4.494.67'
Set-TextValue 'E15' '  -2.48%  '
# Row 16
Set-TextValue 'D16' '3.851.51'
Set-TextValue 'E16' '  -2.23%  '
# Row 17
Set-TextValue 'D17' '68.102.88'
Set-TextValue 'E17' '  -2.52%  '
# Row 18
Set-TextValue 'D18' '17.95'
Set-TextValue 'E18' '  +2.58%  '
# Row 19
Set-TextValue 'D19' '7.31'
Set-TextValue 'E19' '  -3.36%  '
# Row 20
Set-TextValue 'E20' '  -0.61%  '
# Row 21
Set-TextValue 'D21' '10.71'
Set-TextValue 'E21' '  -3.26%  '
# Row 22
Set-TextValue 'D22' '464.90'
Set-TextValue 'E22' '  -7.03%  '
# Row 23
Set-TextValue 'D23' '0.731'
Set-TextValue 'E23' '  -0.81%  '
# Row 24
Set-TextValue 'E24' '  -6.24%  '
# Row 25
Set-TextValue 'D25' '82.90'
Set-TextValue 'E25' '  -3.15%  '
# Row 26
Set-TextValue 'E26' '  -3.40%  '
# Row 27
Set-TextValue 'E27' '  -2.20%  '
# Row 28
Set-TextValue 'E28' '  -0.23%  '
# Row 29
Set-TextValue 'E29' '  -3.53%  '
# Row 30
Set-TextValue 'D30' '2.96'
Set-TextValue 'E30' '  -1.41%  '
# Row 31
Set-TextValue 'D31' '3.998.18'
Set-TextValue 'E31' '  -2.43%  '
# Row 32
Set-TextValue 'E32' '  -3.10%  '
# Row 33
Set-TextValue 'E33' '  -5.01%  '
# Row 34
Set-TextValue 'D34' '31.06'
Set-TextValue 'E34' '  -3.72%  '
# Row 35
Set-TextValue 'D35' '9.41'
Set-TextValue 'E35' '  -1.60%  '
# Row 36
Set-TextValue 'D36' '3.813.72'
Set-TextValue 'E36' '  -2.69%  '
# Row 37
Set-TextValue 'E37' '  -3.57%  '
# Row 38
Set-TextValue 'D38' '3.58'
Set-TextValue 'E38' '  +8.45%  '
# Row 40
Set-TextValue 'E40' '  -0.70%  '
# Row 41
Set-TextValue 'D41' '5.88'
Set-TextValue 'E41' '  -3.79%  '
# Row 42
Set-TextValue 'E42' '  +0.21%  '
# Row 43
Set-TextValue 'D43' '0.313'
Set-TextValue 'E43' '  -3.80%  '
# Row 44
Set-TextValue 'D44' '423.70'
Set-TextValue 'E44' '  -3.07%  '
# Row 45
Set-TextValue 'B45' 'Stacks'
Set-TextValue 'C45' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D45' '1.96'
Set-TextValue 'E45' '  -4.90%  '
# Row 46
Set-TextValue 'B46' 'USDe'
Set-TextValue 'C46' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D46' '1.00'
Set-TextValue 'E46' '  -0.02%  '
# Row 47
Set-TextValue 'B47' 'FLOKI'
Set-TextValue 'C47' 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue 'D47' '0.000293'
Set-TextValue 'E47' '  +5.83%  '
# Row 48
Set-TextValue 'D48' '8.58'
Set-TextValue 'E48' '  -0.30%  '
# Row 49
Set-TextValue 'D49' '46.96'
Set-TextValue 'E49' '  -2.86%  '
# Row 50
Set-TextValue 'D50' '26.29'
Set-TextValue 'E50' '  +2.68%  '
# Row 51
Set-TextValue 'D51' '142.53'
Set-TextValue 'E51' '  -0.52%  '
